$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells being updated so that
# numeric-looking strings (e.g. "1.00", "578.48") are NOT auto-converted
# to actual numbers by Excel, matching the original inline-string text cells.
$priceCells = @("D2","D3","D4","D5","D6","D7","D9","D10","D13","D14","D17","D18","D19","D20","D21","D22","D23","D24","D28","D29","D30","D32","D33","D36","D37","D39","D41","D42","D43","D44","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "61.557.02"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "3.448.87"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "578.48"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "144.85"
$ws.Range("E6").Value = "  +6.78%  "
$ws.Range("D7").Value = "3.449.20"
$ws.Range("E7").Value = "  +2.33%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("D10").Value = "7.62"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  +3.84%  "
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").Value = "4.037.24"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").Value = "28.19"
$ws.Range("E14").Value = "  +8.82%  "
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").Value = "3.452.73"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "61.673.04"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").Value = "6.28"
$ws.Range("E19").Value = "  +8.13%  "
$ws.Range("D20").Value = "14.24"
$ws.Range("E20").Value = "  +4.59%  "
$ws.Range("D21").Value = "9.54"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").Value = "393.70"
$ws.Range("E22").Value = "  +6.00%  "
$ws.Range("D23").Value = "0.564"
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("D24").Value = "73.22"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").Value = "3.592.40"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("D29").Value = "0.178"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("D30").Value = "7.60"
$ws.Range("E30").Value = "  +4.43%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "8.14"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.47"
$ws.Range("E33").Value = "  -8.98%  "
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D36").Value = "24.02"
$ws.Range("E36").Value = "  +3.42%  "
$ws.Range("D37").Value = "3.479.79"
$ws.Range("E38").Value = "  +3.60%  "
$ws.Range("D39").Value = "5.12"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").Value = "167.53"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "0.0782"
$ws.Range("E42").Value = "  +3.48%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "27.84"
$ws.Range("E43").Value = "  +12.14%  "
$ws.Range("D44").Value = "0.804"
$ws.Range("E44").Value = "  +4.43%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +3.13%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "42.28"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "4.48"
$ws.Range("E48").Value = "  +4.34%  "
$ws.Range("D49").Value = "2.600.98"
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("D50").Value = "1.17"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("E51").Value = "  +2.56%  "
